$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 47 and 48: coin/link swapped positions, each refreshed with new scraped price & volume data
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.48%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.95%  "

# Remaining rows: Price (D) and Volume(1h) (E) updates from the latest scrape
$ws.Range("D2").Value = "51.739.50"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "3.088.21"
$ws.Range("E3").Value = "  +3.53%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "388.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.18"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0866"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "3.571.97"
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.79"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.85"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "3.093.22"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.981"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.72"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").Value = "51.810.74"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "0.0₃0971"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.21"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.34"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.64"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("E39").Value = "  +9.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.58"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.67"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.05"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.73%  "
$ws.Range("D49").Value = "2.036.98"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "3.386.99"
$ws.Range("E50").Value = "  +3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.210"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.07%  "
